$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("H2:H27")

# Pick up the Arial/10pt font already used elsewhere in the sheet (style 1)
# via a format-only paste, so no stray/unused font gets created.
$ws.Cells.Item(1, 1).Copy()
$rng.PasteSpecial(-4122)
$rng.NumberFormat = "@"

for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Formula = '="False"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
$ws.Range("H2:H27").Select()
